# Generate Report for Handoff
# Update the "Latest Handoff Date/Datetime" timestamps for the files that were
# just queued for handoff (rows whose status is "Handback transform failed" or
# "Ready for handoff"), across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(4, 6, 7, 8, 9, 10)

# Overview sheet: column D = "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("D$r").Value = "2016-03-19 08:53:55"
}

# zh-cn sheet: column E = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "2016-03-19 08:53:46"
}

# de-de sheet: column E = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "2016-03-19 08:53:55"
}
